# Role-based xlsx report: update/refresh character rows and drop the
# now-redundant trailing row (Katappa's data is folded into row 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Amrendra Bahubali): refresh CharacterID and re-order the
#     Relations / Photos multi-value fields ---
$ws.Range("A2").Value = "64d679773b7e88caeebcc217"
$ws.Range("F2").Value = "64d275752876b66b40a62e28 ; 64d275472876b66b40a62e22"
$ws.Range("G2").Value = "http://localhost:3000/img/amrendra-2.webp ; http://localhost:3000/img/amrendra.jpg"

# --- Row 3: replace "Mahendra Bahubali" with the refreshed "Katappa"
#     record (the old row 4 Katappa record is being superseded) ---
$ws.Range("A3").Value = "64d67a92d56affd78effbea4"
$ws.Range("B3").Value = "Katappa"
$ws.Range("C3").Value = 75
$ws.Range("D3").Value = "male"
$ws.Range("E3").Value = "knight"
$ws.Range("F3").Value = "64d275752876b66b40a62e28 ; 64d275ac2876b66b40a62e2e"
$ws.Range("G3").Value = "http://localhost:3000/img/katappa.jpg"

# --- Row 4 is now obsolete; delete it entirely so the sheet shrinks to
#     A1:G3 ---
$ws.Rows.Item(4).Delete()
